# Insert a new data row for "Vega Modelo de Temuco" (Piña) at row 202,
# pushing the existing rows 202..280 down to 203..281.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(202).Insert()

$ws.Range("A202").Value = 10
$ws.Range("B202").Value = "Vega Modelo de Temuco"
$ws.Range("C202").Value = "La Araucanía"
$ws.Range("D202").Value = 44468
$ws.Range("E202").Value = 9
$ws.Range("F202").Value = "Fruta"
$ws.Range("G202").Value = 100108
$ws.Range("H202").Value = "Tropicales y subtropicales"
$ws.Range("I202").Value = 100108005
$ws.Range("J202").Value = "Piña"
$ws.Range("K202").Value = "Caramelo"
$ws.Range("L202").Value = "Segunda"
$ws.Range("M202").Value = 50
$ws.Range("N202").Value = 18000
$ws.Range("O202").Value = 18000
$ws.Range("P202").Value = 18000
$ws.Range("Q202").Value = "$/caja 14 unidades"
$ws.Range("R202").Value = "Ecuador"
$ws.Range("S202").Value = 1286
$ws.Range("T202").Value = 14
